$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.578.76"
$ws.Range("E2").Value = "  +2.24%  "
$ws.Range("D3").Value = "3.919.59"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "527.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.614"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.13%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -1.52%  "
$ws.Range("E10").Value = "  -1.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000336"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.33"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.20%  "
$ws.Range("D13").Value = "4.544.18"
$ws.Range("E13").Value = "  +0.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.30"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.73%  "
$ws.Range("D15").Value = "3.919.09"
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.55%  "
$ws.Range("E17").Value = "  +8.02%  "
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.83%  "
$ws.Range("D20").Value = "69.598.15"
$ws.Range("E20").Value = "  +2.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "431.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.87%  "
$ws.Range("E22").Value = "  -4.97%  "
$ws.Range("E23").Value = "  -3.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.62"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.50"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "694.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "13.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.94%  "
$ws.Range("E31").Value = "  -2.43%  "
$ws.Range("E32").Value = "  -3.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "67.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +9.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.443"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.42%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.99"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.43%  "
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "0.0₃0861"
$ws.Range("E36").Value = "  -3.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "40.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.150"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0483"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.09"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.63%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.140"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.06%  "
$ws.Range("E47").Value = "  +8.35%  "
$ws.Range("D48").Value = "0.0₆0348"
$ws.Range("E48").Value = "  +4.96%  "
$ws.Range("D49").Value = "2.738.59"
$ws.Range("E49").Value = "  +11.93%  "
$ws.Range("E50").Value = "  -2.75%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "144.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.28%  "
